# This script updates the "想去人数" (number of people interested) column (F)
# on the "展览" (Exhibition) and "全部类型" (All Types) worksheets, reflecting
# refreshed stats output (gh-pages data regeneration).

$wb = $excel.ActiveWorkbook

$exhibitionUpdates = @{
    2 = 1319
    3 = 1199
    5 = 112
    7 = 662
    8 = 107
    9 = 51
    11 = 2392
    12 = 1597
    13 = 1407
    14 = 307
    15 = 236
    16 = 572
    17 = 768
    18 = 54
    19 = 295
    20 = 1091
    24 = 4779
    26 = 319
    27 = 57
    28 = 153
    29 = 135
    32 = 23
    33 = 1030
    34 = 699
    35 = 64
    38 = 378
    39 = 1007
    42 = 156
    43 = 122
}

$allTypesUpdates = @{
    2 = 1319
    5 = 1199
    9 = 112
    11 = 662
    12 = 107
    17 = 2392
    18 = 1597
    19 = 1407
    20 = 307
    21 = 236
    22 = 572
    24 = 768
    25 = 54
    26 = 295
    27 = 1091
    29 = 4779
    31 = 319
    32 = 57
    33 = 153
    34 = 135
    37 = 23
    38 = 1030
    39 = 699
    41 = 378
    42 = 1007
    45 = 156
    46 = 122
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item([int]$row, 6).Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item([int]$row, 6).Value = $allTypesUpdates[$row]
}
